$d = $word.ActiveDocument

# --- Step 1: insert "design " as its own run right before "requirements." ---
# (the first bulleted goal currently reads "...previously described requirements.")
$rng = $d.Content
$rng.Find.Execute("requirements.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$reqStart = $rng.Start

# Typing/inserting plain text merges into the neighbouring run, and toggling
# character formatting to force a run-break leaves a stray empty <w:rPr/>
# behind. Instead, stage the new word as its own run at the very end of the
# document, grab its FormattedText (a distinct run), and paste that at the
# real insertion point - this keeps "design " as its own <w:r> with no
# leftover formatting markers, then we remove the staging text again.
$endPos = $d.Content.End - 1
$endRange = $d.Range($endPos, $endPos)
$endRange.InsertAfter("design ")

$tempRng = $d.Content
$tempRng.Find.Execute("design ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ft = $tempRng.FormattedText

$dest = $d.Range($reqStart, $reqStart)
$dest.FormattedText = $ft

# Remove the temporary staging text again.
$cleanupRng = $d.Content
$cleanupRng.Find.Execute("sections.design ", $true, $false, $false, $false, $false, $true, 1, $false, "sections.", 2) | Out-Null

# --- Step 2: move the "_GoBack" bookmark to sit right before "requirements." ---
# Bookmark names are unique, so adding "_GoBack" again removes it from its
# old location (end of the document) and creates it at the new collapsed
# range - matching Word's own behaviour of re-stamping "_GoBack" at the
# location of the most recent edit.
$rng2 = $d.Content
$rng2.Find.Execute("requirements.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($rng2.Start, $rng2.Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
